# Apply the changes described by the commit:
#  - Rename the sheet "Tabelle1" -> "import this"
#  - Fix the "delivered by" value in row 2 (H2) from "A.S.O.R" to "donor"
#  - Move the window scroll / selection back to the top of the sheet (A1 / H4)
#  - Nudge column S (19) width slightly narrower

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "import this"

# Correct the "delivered by" entry for row 2: it was mistakenly "A.S.O.R",
# should be "donor".
$ws.Range("H2").Value = "donor"

# Reset the view: scroll back to the top-left (A1) and select H4.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H4").Select()

# Slightly narrow column S (column 19).
$ws.Columns.Item(19).ColumnWidth = 15.13
